# This script reproduces the commit:
#   - Two new weekly price records (date serial 44641, "Hass" variety,
#     "Primera" and "Segunda" grades, Provincia de Quillota origin) are
#     inserted at the top of the avocado ("Palta") data block, at rows
#     374-375. All of the existing records that used to occupy rows
#     374-459 are pushed down two rows, ending at row 461, and the
#     worksheet's used range grows from A1:T459 to A1:T461.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 374; this shifts former rows
# 374-459 down to 376-461 (and the row formatting/number format of
# column D is inherited automatically from the surrounding rows).
$ws.Rows("374:375").Insert()

# --- New row 374 --------------------------------------------------
$ws.Range("A374").Value = 7
$ws.Range("B374").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C374").Value = "Ñuble"
$ws.Range("D374").Value = 44641
$ws.Range("E374").Value = 16
$ws.Range("F374").Value = "Fruta"
$ws.Range("G374").Value = 100106
$ws.Range("H374").Value = "Oleaginosos"
$ws.Range("I374").Value = 100106002
$ws.Range("J374").Value = "Palta"
$ws.Range("K374").Value = "Hass"
$ws.Range("L374").Value = "Primera"
$ws.Range("M374").Value = 40
$ws.Range("N374").Value = 3000
$ws.Range("O374").Value = 3000
$ws.Range("P374").Value = 3000
$ws.Range("Q374").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R374").Value = "Provincia de Quillota"
$ws.Range("S374").Value = 3000
$ws.Range("T374").Value = 1

# --- New row 375 --------------------------------------------------
$ws.Range("A375").Value = 7
$ws.Range("B375").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C375").Value = "Ñuble"
$ws.Range("D375").Value = 44641
$ws.Range("E375").Value = 16
$ws.Range("F375").Value = "Fruta"
$ws.Range("G375").Value = 100106
$ws.Range("H375").Value = "Oleaginosos"
$ws.Range("I375").Value = 100106002
$ws.Range("J375").Value = "Palta"
$ws.Range("K375").Value = "Hass"
$ws.Range("L375").Value = "Segunda"
$ws.Range("M375").Value = 80
$ws.Range("N375").Value = 2700
$ws.Range("O375").Value = 2800
$ws.Range("P375").Value = 2750
$ws.Range("Q375").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R375").Value = "Provincia de Quillota"
$ws.Range("S375").Value = 2750
$ws.Range("T375").Value = 1
